$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 72 (pushes existing rows 72..162 down to 73..163)
$ws.Rows.Item(72).Insert()

# Populate the newly inserted row 72 with the new weekly record
$ws.Cells.Item(72, 1).Value = 6
$ws.Cells.Item(72, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(72, 3).Value = "Metropolitana"
$ws.Cells.Item(72, 4).Value = 44477
$ws.Cells.Item(72, 5).Value = 13
$ws.Cells.Item(72, 6).Value = "Fruta"
$ws.Cells.Item(72, 7).Value = 100101
$ws.Cells.Item(72, 8).Value = "Berries"
$ws.Cells.Item(72, 9).Value = 100101001
$ws.Cells.Item(72, 10).Value = "Arándano (blue)"
$ws.Cells.Item(72, 11).Value = "Sin especificar"
$ws.Cells.Item(72, 12).Value = "Especial"
$ws.Cells.Item(72, 13).Value = 750
$ws.Cells.Item(72, 14).Value = 14000
$ws.Cells.Item(72, 15).Value = 14000
$ws.Cells.Item(72, 16).Value = 14000
$ws.Cells.Item(72, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(72, 18).Value = "Perú"
$ws.Cells.Item(72, 19).Value = 7000
$ws.Cells.Item(72, 20).Value = 2
